$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Str2Col (column E) row 3 used to hold the string "b"; change it to the
# number 3 so the object-dtype column has a float/int value sandwiched
# between its string values (a, 3, c, d, e).
$ws.Range("E3").Value = 3

# Move the active selection to E4 (matches the saved view state).
$ws.Range("E4").Select()
